$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$meta.Range("B6").Value = "active"

# Date: updated publish date
$meta.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true
# (Assigning the literal string "true"/"false" via .Value auto-coerces the
# cell to a Boolean in Excel, same as typing it into the grid. Route it
# through a text formula and flatten to a static value via copy/paste-values
# so the result lands back as a genuine text (shared-string) cell.)
$caseSensitive = $meta.Range("B17")
$caseSensitive.Formula = "=""true"""
$caseSensitive.Copy()
$caseSensitive.PasteSpecial(-4163)
